$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same table of events and
# need their "想去人数" (want-to-go count) values updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8800
    $ws.Range("F4").Value = 425
    $ws.Range("F5").Value = 239
}
